$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.10769139134775
$ws.Range("G2").Value = 97520.2281613775
$ws.Range("H2").Value = -194874.456322755
$ws.Range("K2").Value = -194102.805874818

$ws.Range("F3").Value = 0.108080375720078
$ws.Range("G3").Value = 96615.7110242586
$ws.Range("H3").Value = -193195.422048517
$ws.Range("I3").Value = 1679.03427423778
$ws.Range("K3").Value = -193028.076168242

$ws.Range("F4").Value = 0.10771065056736
$ws.Range("G4").Value = 96249.9594431749
$ws.Range("H4").Value = -192455.91888635
$ws.Range("I4").Value = 2418.53743640525
$ws.Range("K4").Value = -192251.38503268

$ws.Range("F5").Value = 0.108351380927765
$ws.Range("G5").Value = 96220.8707511761
$ws.Range("H5").Value = -192397.741502352
$ws.Range("I5").Value = 2476.71482040279
$ws.Range("K5").Value = -192193.207648682

$ws.Range("I6").Value = 2835.4564113231

$ws.Range("F7").Value = 0.108436433145375
$ws.Range("G7").Value = 95845.7489746352
$ws.Range("H7").Value = -191679.49794927
$ws.Range("I7").Value = 3194.95837348461
$ws.Range("K7").Value = -191623.715989179

$ws.Range("F8").Value = 0.108084000859814
$ws.Range("G8").Value = 95419.8018318762
$ws.Range("H8").Value = -190825.603663752
$ws.Range("I8").Value = 4048.85265900267
$ws.Range("K8").Value = -190760.524710312
